# Applies the NPA-NNAC-08E (NTPJ) revision:
#  1. "ATRIBUICOES" heading -> "RESPONSABILIDADE POR FUNCAO"
#  2. The explanatory paragraph below it is reworded.
#  3. The "armazenamento e arquivo de Desenhos e ADT" bullet becomes
#     "armazenamento de Desenhos".
#  4. Two new numbered items are appended describing the
#     inter-sector responsibility.

$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- 1. Section heading -------------------------------------------------
$d.Content.Find.Execute(
    "ATRIBUIÇÕES", $false, $false, $false, $false, $false,
    $true, 1, $false, "RESPONSABILIDADE POR FUNÇÃO", 2) | Out-Null

# --- 2. Locate the explanatory paragraph + the "armazenamento" item -----
# (re-located by content since indices shift as we edit)
$introPara = $null
$desenhosPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Cabe a Chefe da*garantir o cumprimento dos seguintes processos*") {
        $introPara = $p
    }
    if ($t -like "*rmazenamento e arquivo de Desenhos e ADT*") {
        $desenhosPara = $p
    }
}

# --- 2a. Reword the explanatory paragraph, preserving its own formatting
$introXml = "<w:p $wNs>" +
    "<w:pPr><w:spacing w:before=`"120`"/><w:rPr><w:sz w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr><w:t xml:space=`"preserve`">As </w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr><w:t>responsabilidades por função são delineadas sinteticamente no Regimento Interno do CELOG e as atividades relacionados aos seus cumprimentos estão contempladas n</w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr><w:t>os seguintes processos:</w:t></w:r>" +
    "</w:p>"
$introPara.Range.InsertXML($introXml)

# --- 2b. Split the last run of the "armazenamento" bullet ---------------
$desenhosXml = "<w:p $wNs>" +
    "<w:pPr><w:keepNext/><w:keepLines/><w:widowControl w:val=`"0`"/><w:numPr><w:ilvl w:val=`"2`"/><w:numId w:val=`"1`"/></w:numPr><w:spacing w:before=`"120`"/><w:jc w:val=`"both`"/><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr><w:t>PLOG</w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr><w:t>0</w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr><w:t>0</w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr><w:t>18</w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr><w:t xml:space=`"preserve`"> – Cadastro e </w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr><w:t>a</w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr><w:t xml:space=`"preserve`">rmazenamento </w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr><w:t>de</w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr><w:t xml:space=`"preserve`"> Desenhos</w:t></w:r>" +
    "</w:p>"
$desenhosPara.Range.InsertXML($desenhosXml)

# --- 3. Insert the two new numbered items after the "Desenhos" bullet ---
# Re-find the bullet (it was just reseated by InsertXML above) and work
# from the start of the blank paragraph that immediately follows it; that
# blank paragraph's own formatting is reproduced as the trailing 4th
# fragment below so it survives untouched.
$desenhosPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Desenhos*") {
        $desenhosPara = $p
    }
}
$nextPara = $desenhosPara.Next()
$insertRange = $nextPara.Range
$insertRange.Collapse(1)

$blankPPr = "<w:pPr><w:keepNext/><w:keepLines/><w:widowControl w:val=`"0`"/><w:spacing w:before=`"120`"/><w:jc w:val=`"both`"/><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr></w:pPr>"

$newXml =
    "<w:p $wNs>$blankPPr</w:p>" +
    "<w:p $wNs><w:pPr><w:keepNext/><w:keepLines/><w:widowControl w:val=`"0`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr><w:spacing w:before=`"120`"/><w:jc w:val=`"both`"/><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:u w:val=`"single`"/><w:lang w:val=`"pt-PT`"/></w:rPr><w:t>RESPONSABILIDADE NO INTER-RELACIONAMENTO ENTRE OS SETORES</w:t></w:r></w:p>" +
    "<w:p $wNs><w:pPr><w:keepNext/><w:keepLines/><w:widowControl w:val=`"0`"/><w:numPr><w:ilvl w:val=`"2`"/><w:numId w:val=`"1`"/></w:numPr><w:spacing w:before=`"120`"/><w:jc w:val=`"both`"/><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr><w:t>As responsa</w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:lang w:val=`"pt-PT`"/></w:rPr><w:t>bilidades no inter-relacionamento entre setores são apresentadas detalhadamente nos PLOG relacionados no item 2.2 desta NPA.</w:t></w:r></w:p>" +
    "<w:p $wNs>$blankPPr</w:p>"

$insertRange.InsertXML($newXml)

# --- 4. Bump the cached PAGE-number field result in the header ----------
# ("fl. 2" -> "fl. 3": the document now spans one more page). The engine
# does not re-run full pagination, so the cached field result is patched
# directly (mirrors what Word itself stores as the field's last-computed
# result).
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$hdrRange = $hdr.Range
for ($i = 1; $i -le $hdrRange.Characters.Count; $i++) {
    $ch = $hdrRange.Characters.Item($i)
    if ($ch.Text -eq "2") {
        $ch.Text = "3"
        break
    }
}

Write-Output "Done"
